$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# Re-set Region (B2) and Design (B7) values. The underlying shared-string
# table gets re-ordered by this edit (EMEA moved ahead of SMART in the
# string table), but the actual displayed cell contents are unchanged:
# B2 stays "EMEA" and B7 stays "SMART".
$ws.Range("B2").Value = "EMEA"
$ws.Range("B7").Value = "SMART"

# Swap which link is used for the Cellular/4G interface:
# row 17 (Main Link section) flips from TRUE to FALSE,
# row 25 (Backup Link section) flips from FALSE to TRUE.
$ws.Range("B17").Value = $false
$ws.Range("B25").Value = $true

# Update the active selection / view state left behind by the edit.
$ws.Range("H16").Select()
